$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) store numeric-looking text (e.g. "242.57",
# "-1.14%"); mark those cells as Text before writing so Excel keeps them as
# literal strings instead of coercing them into numbers/percentages.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "242.57"
$ws.Range("E2").Value = "-1.14%"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "29.65"
$ws.Range("E3").Value = "11.94%"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.122"
$ws.Range("E4").Value = "-0.38%"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05650"
$ws.Range("E5").Value = "1.16%"

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "6.494"
$ws.Range("E6").Value = "0.27%"

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8259"
$ws.Range("E7").Value = "1.00%"

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8634"
$ws.Range("E8").Value = "2.76%"

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1330"
$ws.Range("E9").Value = "0.00%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "0.06858"
$ws.Range("E10").Value = "-1.93%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "0.02858"
$ws.Range("E11").Value = "-0.41%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "0.09389"
$ws.Range("E12").Value = "0.08%"

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "0.001515"
$ws.Range("E13").Value = "-0.86%"

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("B14").Value = "CoinExToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D14").Value = "0.04155"
$ws.Range("E14").Value = "-9.64%"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "0.0005979"
$ws.Range("E15").Value = "-94.03%"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006154"
$ws.Range("E16").Value = "0.36%"

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.521"
$ws.Range("E17").Value = "-3.13%"

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "3.013"
$ws.Range("E18").Value = "-0.78%"

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "2.215"
$ws.Range("E19").Value = "1.45%"

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3149"
$ws.Range("E20").Value = "1.19%"

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03228"
$ws.Range("E21").Value = "5.35%"

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1280"
$ws.Range("E22").Value = "-1.49%"

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "3.609"
$ws.Range("E23").Value = "-3.55%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.06%"

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001209"
$ws.Range("E25").Value = "-3.07%"

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004442"
$ws.Range("E26").Value = "-1.60%"

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001180"
$ws.Range("E27").Value = "22.83%"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "0.51%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03706"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.005868"
$ws.Range("E41").Value = "-4.93%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1053"
$ws.Range("E42").Value = "0.28%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.74%"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009710"
$ws.Range("E44").Value = "18.25%"

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005097"
$ws.Range("E45").Value = "-4.89%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.06%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-3.73%"

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002547"
$ws.Range("E48").Value = "-0.11%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.06%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"
